$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-02 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-03 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("282×8=2256", $true, $false, $false, $false, $false, $true, 1, $false, "966×3=2898", 2) | Out-Null
$d.Content.Find.Execute("690×8=5520", $true, $false, $false, $false, $false, $true, 1, $false, "648×9=5832", 2) | Out-Null
$d.Content.Find.Execute("677×2=1354", $true, $false, $false, $false, $false, $true, 1, $false, "965×5=4825", 2) | Out-Null
$d.Content.Find.Execute("747×6=4482", $true, $false, $false, $false, $false, $true, 1, $false, "932×8=7456", 2) | Out-Null
$d.Content.Find.Execute("999×2=1998", $true, $false, $false, $false, $false, $true, 1, $false, "672×6=4032", 2) | Out-Null
$d.Content.Find.Execute("894×8=7152", $true, $false, $false, $false, $false, $true, 1, $false, "344×3=1032", 2) | Out-Null
$d.Content.Find.Execute("151×9=1359", $true, $false, $false, $false, $false, $true, 1, $false, "381×9=3429", 2) | Out-Null
$d.Content.Find.Execute("864×6=5184", $true, $false, $false, $false, $false, $true, 1, $false, "209×8=1672", 2) | Out-Null
$d.Content.Find.Execute("993×7=6951", $true, $false, $false, $false, $false, $true, 1, $false, "142×9=1278", 2) | Out-Null
$d.Content.Find.Execute("974×4=3896", $true, $false, $false, $false, $false, $true, 1, $false, "884×6=5304", 2) | Out-Null
$d.Content.Find.Execute("979×5=4895", $true, $false, $false, $false, $false, $true, 1, $false, "474×8=3792", 2) | Out-Null
$d.Content.Find.Execute("409×4=1636", $true, $false, $false, $false, $false, $true, 1, $false, "514×4=2056", 2) | Out-Null
$d.Content.Find.Execute("548×9=4932", $true, $false, $false, $false, $false, $true, 1, $false, "273×6=1638", 2) | Out-Null
$d.Content.Find.Execute("256×6=1536", $true, $false, $false, $false, $false, $true, 1, $false, "665×8=5320", 2) | Out-Null
$d.Content.Find.Execute("592×9=5328", $true, $false, $false, $false, $false, $true, 1, $false, "880×9=7920", 2) | Out-Null
$d.Content.Find.Execute("581×7=4067", $true, $false, $false, $false, $false, $true, 1, $false, "706×9=6354", 2) | Out-Null
$d.Content.Find.Execute("490×4=1960", $true, $false, $false, $false, $false, $true, 1, $false, "717×3=2151", 2) | Out-Null
$d.Content.Find.Execute("968×9=8712", $true, $false, $false, $false, $false, $true, 1, $false, "528×3=1584", 2) | Out-Null
$d.Content.Find.Execute("866×2=1732", $true, $false, $false, $false, $false, $true, 1, $false, "119×2=238", 2) | Out-Null
$d.Content.Find.Execute("375×3=1125", $true, $false, $false, $false, $false, $true, 1, $false, "768×6=4608", 2) | Out-Null
$d.Content.Find.Execute("245×6=1470", $true, $false, $false, $false, $false, $true, 1, $false, "236×9=2124", 2) | Out-Null
$d.Content.Find.Execute("663×4=2652", $true, $false, $false, $false, $false, $true, 1, $false, "127×6=762", 2) | Out-Null
$d.Content.Find.Execute("748×9=6732", $true, $false, $false, $false, $false, $true, 1, $false, "785×6=4710", 2) | Out-Null
$d.Content.Find.Execute("364×9=3276", $true, $false, $false, $false, $false, $true, 1, $false, "650×8=5200", 2) | Out-Null
$d.Content.Find.Execute("862×4=3448", $true, $false, $false, $false, $false, $true, 1, $false, "768×5=3840", 2) | Out-Null
